$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").Value = "group"
$ws.Range("D2:D5").Value = "ETS_CO2"
$ws.Range("D2:D4").Style = "Normal"
$ws.Range("D2:D5").Select()
